# daily auto push: 2026-01-17 13:37 UTC
# Insert a new data row for 2026/01/17 (Sat) 19:00 slot before the existing
# row 643, shifting the 2026/12/29 .. 2027/01/05 block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything currently at/after row 643 down by one row.
$ws.Range("A643").EntireRow.Insert()

# Column A holds a date formatted as literal text (e.g. "2026/01/17"), not
# a real date serial. Force the cell to Text format before assigning so
# Excel doesn't auto-convert the string into a date number, then drop the
# number format again so the cell matches its plain, unstyled neighbours.
$dateCell = $ws.Cells.Item(643, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/17"
$dateCell.ClearFormats()

$ws.Cells.Item(643, 2).Value = "土"
$ws.Cells.Item(643, 3).Value = 19
$ws.Cells.Item(643, 4).Value = 23
